$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new header cell
$ws.Range("C1").Value = "Ник пользователя"

# Apply the same style (font + fill) used for the header row to the whole A1:C1 range
$ws.Range("A1:C1").Font.Name = "Calibri"
$ws.Range("A1:C1").Font.Size = 11
$ws.Range("A1:C1").Interior.Pattern = -4124
$ws.Range("A1:C1").Interior.ThemeColor = 9
$ws.Range("A1:C1").Interior.TintAndShade = 0.39997558519241921

# Column widths: columns A-C should all be 30.7109375 wide
$ws.Range("A1:C1").ColumnWidth = 30.7109375

# Update selection to A2
$ws.Range("A2").Select()
